# Apply forecast summary correction: insert Week_Start_Date column,
# shorten week labels, and mark is_holiday_week as boolean.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before column B (ASIN) to make room for Week_Start_Date
$ws.Columns.Item(2).Insert()

# Header row
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# Week start dates (Sundays), 16 weeks starting 2025-01-05
$weekStarts = @(
    "2025-01-05", "2025-01-12", "2025-01-19", "2025-01-26",
    "2025-02-02", "2025-02-09", "2025-02-16", "2025-02-23",
    "2025-03-02", "2025-03-09", "2025-03-16", "2025-03-23",
    "2025-03-30", "2025-04-06", "2025-04-13", "2025-04-20"
)

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2

    # Shorten week label, e.g. "W01" -> "W1"
    $ws.Cells.Item($row, 1).Value = "W" + ($i + 1)

    # New Week_Start_Date column stored as literal text (not a date serial)
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $weekStarts[$i]

    # is_holiday_week (now column J) should be boolean FALSE
    $ws.Cells.Item($row, 10).Value = $false
}
